$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1443
$ws.Range("F3").Value = 108
$ws.Range("F4").Value = 2110
$ws.Range("F5").Value = 6960
$ws.Range("F6").Value = 551
$ws.Range("F7").Value = 1059
$ws.Range("F8").Value = 48
$ws.Range("F9").Value = 4620
$ws.Range("F10").Value = 6837
$ws.Range("F11").Value = 13
$ws.Range("F12").Value = 234
$ws.Range("F13").Value = 1410
$ws.Range("F14").Value = 825
$ws.Range("F15").Value = 124
$ws.Range("F17").Value = 36
$ws.Range("F20").Value = 138
$ws.Range("F22").Value = 193
$ws.Range("F24").Value = 1080
$ws.Range("F26").Value = 39
$ws.Range("F28").Value = 31
$ws.Range("F29").Value = 123
$ws.Range("F32").Value = 112
$ws.Range("G32").Value = 99
$ws.Range("F33").Value = 9
$ws.Range("F35").Value = 3
$ws.Range("F40").Value = 50
$ws.Range("F43").Value = 1185
$ws.Range("F45").Value = 67
$ws.Range("F46").Value = 120

$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 21
$ws.Range("F3").Value = 21
$ws.Range("F13").Value = 20
$ws.Range("F18").Value = 550
$ws.Range("F24").Value = 136
$ws.Range("F31").Value = 819
$ws.Range("F36").Value = 95
$ws.Range("F40").Value = 130

$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 847
$ws.Range("F6").Value = 618
$ws.Range("F8").Value = 1355
$ws.Range("F9").Value = 2131

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 21
$ws.Range("F4").Value = 1443
$ws.Range("F6").Value = 847
$ws.Range("F7").Value = 108
$ws.Range("F8").Value = 618
$ws.Range("F9").Value = 618
$ws.Range("F11").Value = 6960
$ws.Range("F12").Value = 551
$ws.Range("F13").Value = 48
$ws.Range("F14").Value = 4620
$ws.Range("F15").Value = 6837
$ws.Range("F16").Value = 13
$ws.Range("F17").Value = 234
$ws.Range("F18").Value = 1410
$ws.Range("F20").Value = 825
$ws.Range("F21").Value = 124
$ws.Range("F22").Value = 1355
$ws.Range("F23").Value = 2131
$ws.Range("F25").Value = 36
$ws.Range("F27").Value = 138
$ws.Range("F28").Value = 193
$ws.Range("F29").Value = 1080
$ws.Range("F33").Value = 123
$ws.Range("F35").Value = 112
$ws.Range("G35").Value = 99
$ws.Range("F37").Value = 819
$ws.Range("F43").Value = 50
$ws.Range("F44").Value = 95
$ws.Range("F49").Value = 120
